$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the new survival probability: update the resulting CE values
$ws.Range("B2").Value = 15052.56464127515
$ws.Range("C2").Value = 824814.8374961047

$ws.Range("B3").Value = 22348.25969763185
$ws.Range("C3").Value = 1224587.080681118

$ws.Range("B4").Value = 34777.10901857738
$ws.Range("C4").Value = 1905633.771210448
